$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update nombre_aides (col C) and montant_total (col D) for rows with new data
# as part of the 2020-12-28 data refresh for Fonds de solidarite volet 1

$ws.Cells.Item(2, 3).Value = 56109
$ws.Cells.Item(2, 4).Value = 114229261
$ws.Cells.Item(3, 3).Value = 136160
$ws.Cells.Item(3, 4).Value = 319282363
$ws.Cells.Item(4, 3).Value = 49410
$ws.Cells.Item(4, 4).Value = 143036325
$ws.Cells.Item(5, 3).Value = 15610
$ws.Cells.Item(5, 4).Value = 52500774
$ws.Cells.Item(6, 3).Value = 5740
$ws.Cells.Item(6, 4).Value = 26062899
$ws.Cells.Item(7, 3).Value = 1133
$ws.Cells.Item(7, 4).Value = 6592279
$ws.Cells.Item(8, 3).Value = 68
$ws.Cells.Item(8, 4).Value = 469315
$ws.Cells.Item(12, 3).Value = 58301
$ws.Cells.Item(12, 4).Value = 92669294
$ws.Cells.Item(13, 3).Value = 14239
$ws.Cells.Item(13, 4).Value = 28639924
$ws.Cells.Item(14, 3).Value = 38325
$ws.Cells.Item(14, 4).Value = 88336075
$ws.Cells.Item(15, 3).Value = 12733
$ws.Cells.Item(15, 4).Value = 35159617
$ws.Cells.Item(16, 3).Value = 3665
$ws.Cells.Item(16, 4).Value = 11229437
$ws.Cells.Item(17, 3).Value = 1195
$ws.Cells.Item(17, 4).Value = 5128241
$ws.Cells.Item(18, 3).Value = 256
$ws.Cells.Item(18, 4).Value = 1376157
$ws.Cells.Item(20, 3).Value = 14306
$ws.Cells.Item(20, 4).Value = 22167461
$ws.Cells.Item(21, 3).Value = 20012
$ws.Cells.Item(21, 4).Value = 42247731
$ws.Cells.Item(22, 3).Value = 47440
$ws.Cells.Item(22, 4).Value = 113864701
$ws.Cells.Item(23, 3).Value = 16436
$ws.Cells.Item(23, 4).Value = 47323769
$ws.Cells.Item(24, 3).Value = 4855
$ws.Cells.Item(24, 4).Value = 15656106
$ws.Cells.Item(25, 3).Value = 1552
$ws.Cells.Item(25, 4).Value = 6450394
$ws.Cells.Item(26, 3).Value = 249
$ws.Cells.Item(26, 4).Value = 1361338
$ws.Cells.Item(28, 3).Value = 16014
$ws.Cells.Item(28, 4).Value = 24680790
$ws.Cells.Item(29, 3).Value = 11313
$ws.Cells.Item(29, 4).Value = 23070337
$ws.Cells.Item(30, 3).Value = 32781
$ws.Cells.Item(30, 4).Value = 73916368
$ws.Cells.Item(31, 3).Value = 11819
$ws.Cells.Item(31, 4).Value = 31837648
$ws.Cells.Item(32, 3).Value = 3244
$ws.Cells.Item(32, 4).Value = 9679015
$ws.Cells.Item(33, 3).Value = 1017
$ws.Cells.Item(33, 4).Value = 4291496
$ws.Cells.Item(36, 3).Value = 11581
$ws.Cells.Item(36, 4).Value = 17988706
$ws.Cells.Item(37, 3).Value = 5006
$ws.Cells.Item(37, 4).Value = 10978059
$ws.Cells.Item(38, 3).Value = 11762
$ws.Cells.Item(38, 4).Value = 27283967
$ws.Cells.Item(39, 3).Value = 4865
$ws.Cells.Item(39, 4).Value = 13817360
$ws.Cells.Item(40, 3).Value = 1352
$ws.Cells.Item(40, 4).Value = 4397400
$ws.Cells.Item(41, 3).Value = 431
$ws.Cells.Item(41, 4).Value = 2056184
$ws.Cells.Item(44, 3).Value = 3560
$ws.Cells.Item(44, 4).Value = 5477171
$ws.Cells.Item(45, 3).Value = 25629
$ws.Cells.Item(45, 4).Value = 52882761
$ws.Cells.Item(46, 3).Value = 76113
$ws.Cells.Item(46, 4).Value = 179356347
$ws.Cells.Item(47, 3).Value = 29263
$ws.Cells.Item(47, 4).Value = 81407767
$ws.Cells.Item(48, 3).Value = 9521
$ws.Cells.Item(48, 4).Value = 28938770
$ws.Cells.Item(49, 3).Value = 3278
$ws.Cells.Item(49, 4).Value = 13163538
$ws.Cells.Item(50, 3).Value = 562
$ws.Cells.Item(50, 4).Value = 3170046
$ws.Cells.Item(53, 3).Value = 26204
$ws.Cells.Item(53, 4).Value = 47955009
$ws.Cells.Item(55, 3).Value = 9020
$ws.Cells.Item(55, 4).Value = 14917774
$ws.Cells.Item(56, 3).Value = 3031
$ws.Cells.Item(56, 4).Value = 5345477
$ws.Cells.Item(57, 3).Value = 991
$ws.Cells.Item(57, 4).Value = 1936183
$ws.Cells.Item(58, 3).Value = 300
$ws.Cells.Item(58, 4).Value = 667437
$ws.Cells.Item(61, 3).Value = 9216
$ws.Cells.Item(61, 4).Value = 13668842
$ws.Cells.Item(62, 3).Value = 1823
$ws.Cells.Item(62, 4).Value = 3973977
$ws.Cells.Item(63, 3).Value = 4311
$ws.Cells.Item(63, 4).Value = 9343978
$ws.Cells.Item(64, 3).Value = 1723
$ws.Cells.Item(64, 4).Value = 3868959
$ws.Cells.Item(68, 3).Value = 2827
$ws.Cells.Item(68, 4).Value = 5590039
$ws.Cells.Item(69, 3).Value = 22841
$ws.Cells.Item(69, 4).Value = 45114762
$ws.Cells.Item(70, 3).Value = 66280
$ws.Cells.Item(70, 4).Value = 151028229
$ws.Cells.Item(71, 3).Value = 24215
$ws.Cells.Item(71, 4).Value = 67034286
$ws.Cells.Item(72, 3).Value = 7573
$ws.Cells.Item(72, 4).Value = 22850275
$ws.Cells.Item(73, 3).Value = 2450
$ws.Cells.Item(73, 4).Value = 9896879
$ws.Cells.Item(74, 3).Value = 484
$ws.Cells.Item(74, 4).Value = 2673913
$ws.Cells.Item(78, 3).Value = 21190
$ws.Cells.Item(78, 4).Value = 32519668
$ws.Cells.Item(79, 3).Value = 83270
$ws.Cells.Item(79, 4).Value = 170790139
$ws.Cells.Item(80, 3).Value = 226310
$ws.Cells.Item(80, 4).Value = 509463899
$ws.Cells.Item(81, 3).Value = 102133
$ws.Cells.Item(81, 4).Value = 285693395
$ws.Cells.Item(82, 3).Value = 37236
$ws.Cells.Item(82, 4).Value = 124951614
$ws.Cells.Item(83, 3).Value = 13683
$ws.Cells.Item(83, 4).Value = 61649660
$ws.Cells.Item(84, 3).Value = 2656
$ws.Cells.Item(84, 4).Value = 17016509
$ws.Cells.Item(85, 3).Value = 153
$ws.Cells.Item(85, 4).Value = 888967
$ws.Cells.Item(90, 3).Value = 79086
$ws.Cells.Item(90, 4).Value = 125555569
$ws.Cells.Item(91, 3).Value = 5624
$ws.Cells.Item(91, 4).Value = 8760651
$ws.Cells.Item(92, 3).Value = 13553
$ws.Cells.Item(92, 4).Value = 21429147
$ws.Cells.Item(93, 3).Value = 4353
$ws.Cells.Item(93, 4).Value = 7044628
$ws.Cells.Item(94, 3).Value = 1540
$ws.Cells.Item(94, 4).Value = 2569024
$ws.Cells.Item(98, 3).Value = 6358
$ws.Cells.Item(98, 4).Value = 8755989
$ws.Cells.Item(99, 3).Value = 2307
$ws.Cells.Item(99, 4).Value = 4273860
$ws.Cells.Item(100, 3).Value = 7312
$ws.Cells.Item(100, 4).Value = 14527982
$ws.Cells.Item(102, 3).Value = 969
$ws.Cells.Item(102, 4).Value = 2305268
$ws.Cells.Item(106, 3).Value = 4895
$ws.Cells.Item(106, 4).Value = 7120564
$ws.Cells.Item(107, 3).Value = 1034
$ws.Cells.Item(107, 4).Value = 2334197
$ws.Cells.Item(113, 3).Value = 16377
$ws.Cells.Item(113, 4).Value = 34364941
$ws.Cells.Item(114, 3).Value = 43105
$ws.Cells.Item(114, 4).Value = 100895497
$ws.Cells.Item(115, 3).Value = 15278
$ws.Cells.Item(115, 4).Value = 42633772
$ws.Cells.Item(116, 3).Value = 4709
$ws.Cells.Item(116, 4).Value = 14843816
$ws.Cells.Item(117, 3).Value = 1459
$ws.Cells.Item(117, 4).Value = 6173650
$ws.Cells.Item(118, 3).Value = 302
$ws.Cells.Item(118, 4).Value = 1705408
$ws.Cells.Item(122, 3).Value = 13502
$ws.Cells.Item(122, 4).Value = 20681950
$ws.Cells.Item(123, 3).Value = 44622
$ws.Cells.Item(123, 4).Value = 90323194
$ws.Cells.Item(124, 3).Value = 95513
$ws.Cells.Item(124, 4).Value = 214702796
$ws.Cells.Item(125, 3).Value = 32277
$ws.Cells.Item(125, 4).Value = 86577082
$ws.Cells.Item(126, 3).Value = 10104
$ws.Cells.Item(126, 4).Value = 30630059
$ws.Cells.Item(127, 3).Value = 3193
$ws.Cells.Item(127, 4).Value = 13073837
$ws.Cells.Item(128, 3).Value = 646
$ws.Cells.Item(128, 4).Value = 3491290
$ws.Cells.Item(132, 3).Value = 34984
$ws.Cells.Item(132, 4).Value = 53504362
$ws.Cells.Item(133, 3).Value = 53709
$ws.Cells.Item(133, 4).Value = 110087932
$ws.Cells.Item(134, 3).Value = 112438
$ws.Cells.Item(134, 4).Value = 250888944
$ws.Cells.Item(135, 3).Value = 36412
$ws.Cells.Item(135, 4).Value = 100005175
$ws.Cells.Item(136, 3).Value = 10775
$ws.Cells.Item(136, 4).Value = 33222496
$ws.Cells.Item(137, 3).Value = 3376
$ws.Cells.Item(137, 4).Value = 13884927
$ws.Cells.Item(138, 3).Value = 539
$ws.Cells.Item(138, 4).Value = 2952989
$ws.Cells.Item(142, 3).Value = 43707
$ws.Cells.Item(142, 4).Value = 65498362
$ws.Cells.Item(143, 3).Value = 19603
$ws.Cells.Item(143, 4).Value = 40173370
$ws.Cells.Item(144, 3).Value = 47678
$ws.Cells.Item(144, 4).Value = 111973086
$ws.Cells.Item(145, 3).Value = 17911
$ws.Cells.Item(145, 4).Value = 50031605
$ws.Cells.Item(146, 3).Value = 5150
$ws.Cells.Item(146, 4).Value = 15841404
$ws.Cells.Item(147, 3).Value = 1530
$ws.Cells.Item(147, 4).Value = 6459569
$ws.Cells.Item(148, 3).Value = 342
$ws.Cells.Item(148, 4).Value = 1982434
$ws.Cells.Item(152, 3).Value = 14779
$ws.Cells.Item(152, 4).Value = 22909362
$ws.Cells.Item(153, 3).Value = 53216
$ws.Cells.Item(153, 4).Value = 110057087
$ws.Cells.Item(154, 3).Value = 123651
$ws.Cells.Item(154, 4).Value = 285000278
$ws.Cells.Item(155, 3).Value = 39333
$ws.Cells.Item(155, 4).Value = 112744656
$ws.Cells.Item(156, 3).Value = 11749
$ws.Cells.Item(156, 4).Value = 39419081
$ws.Cells.Item(157, 3).Value = 4208
$ws.Cells.Item(157, 4).Value = 18813979
$ws.Cells.Item(158, 3).Value = 852
$ws.Cells.Item(158, 4).Value = 5184370
$ws.Cells.Item(159, 3).Value = 52
$ws.Cells.Item(159, 4).Value = 246597
$ws.Cells.Item(160, 3).Value = 40986
$ws.Cells.Item(160, 4).Value = 63650757
